$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (bold, bordered, centered) onto the new
# I1:J1 header cells so they match the rest of the header row formatting.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New header labels for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for the new columns I0 (I) and IF (J), rows 2-35
$data = @(
    @(8, 8),
    @(8, 8),
    @(11, 11),
    @(9, 9),
    @(8, 9),
    @(8, 8),
    @(8, 9),
    @(5, 6),
    @(5, 6),
    @(6, 6),
    @(7, 7),
    @(5, 6),
    @(7, 8),
    @(7, 8),
    @(7, 8),
    @(7, 8),
    @(4, 6),
    @(6, 7),
    @(5, 5),
    @(8, 8),
    @(6, 6),
    @(7, 7),
    @(5, 5),
    @(5, 6),
    @(9, 9),
    @(6, 8),
    @(2, 3),
    @(7, 7),
    @(4, 5),
    @(3, 4),
    @(8, 8),
    @(2, 3),
    @(3, 4),
    @(4, 4)
)

for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $data[$idx][0]
    $ws.Cells.Item($row, 10).Value = $data[$idx][1]
}
